$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$s.Shapes.Item("Group 8").Delete()
